$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows of data (rows 74 and 75)
$ws.Cells.Item(74, 1).Value = 45031
$ws.Cells.Item(74, 2).Value = 0.8175810185185185
$ws.Cells.Item(74, 3).Value = 75573
$ws.Cells.Item(74, 4).Value = 1430

$ws.Cells.Item(75, 1).Value = 41382
$ws.Cells.Item(75, 2).Value = 0.65328703703703705
$ws.Cells.Item(75, 3).Value = 75573
$ws.Cells.Item(75, 4).Value = 1430

# Match the date/time formatting used by the existing rows (col A = date, col B = time)
# by copying the style from the row above, so the existing numFmt styles are reused
# instead of creating duplicate style entries.
$ws.Range("A73:B73").Copy() | Out-Null
$ws.Range("A74:B75").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update selection to mimic the recorded state after appending rows
$ws.Range("A76").Select() | Out-Null
